$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) figures on both the "展览" and
# "全部类型" sheets, which carry duplicate copies of the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1890
    $ws.Range("F3").Value = 356
    $ws.Range("F4").Value = 1153
    $ws.Range("F5").Value = 1180
    $ws.Range("F7").Value = 5983
}
